$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: event cancelled, attendee count and price change
    $ws.Range("C2").Value = "合肥·星芒1.5动漫嘉年华（取消）"
    $ws.Range("F2").Value = 1391
    $ws.Range("G2").Value = "不可售"

    # Row 3: attendee count change
    $ws.Range("F3").Value = 2251

    # Row 4: attendee count change
    $ws.Range("F4").Value = 368

    # Row 5: attendee count change
    $ws.Range("F5").Value = 79

    # Row 6: attendee count change
    $ws.Range("F6").Value = 6432

    # Row 7: attendee count change
    $ws.Range("F7").Value = 304
}
